$wb = $excel.ActiveWorkbook

# The new values to place into C11:I11 on each worksheet
$values = @(
    0.31228414330799126,
    0.2127161482774973,
    0.6832700485824148,
    -0.17699999999999982,
    1.5829618029997903,
    16.12947350163202,
    -0.134442166553219
)

foreach ($ws in $wb.Worksheets) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = 3 + $i   # Column C is index 3
        $ws.Cells.Item(11, $col).Value = $values[$i]
    }
}
